$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect so the cells can be edited,
# then re-protect it afterwards to restore the sheet's protected state.
$ws.Unprotect()

# Update the "as of" date in the disclosure text from 2021-05-18 to 2021-05-19
# (only the two digits that differ are touched, leaving the rest of the text intact)
$disclosureCell = $ws.Range("A9")
$dateChars = $disclosureCell.Characters(114, 10)
$dateChars.Text = "2021-05-19"

# Update the Weight (D) and Percent Change (E) figures for rows 2-6
$ws.Range("D2").Value = 0.2589867191116467
$ws.Range("E2").Value = -0.01527936145952113

$ws.Range("D3").Value = 0.2550329349032841
$ws.Range("E3").Value = -0.006146445750935481

$ws.Range("D4").Value = 0.2435951465581892
$ws.Range("E4").Value = -0.005637088152395719

$ws.Range("D5").Value = 0.24238519942688
$ws.Range("E5").Value = 0.001585204755614189

$ws.Range("E6").Value = -0.006513634937429624

$ws.Protect()
